$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# are forced to stay text via a temporary "@" (Text) number format, then the
# format is reset back to the default "Normal" style so no stray formatting is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.531'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.489'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0803'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.38'
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.787'
$ws.Range("D17").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.61'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0746'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.05'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.14'
$ws.Range("D42").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.15'
$ws.Range("D51").Style = "Normal"

# Remaining cells already contain non-numeric-looking text (percentages with
# surrounding spaces/%, or multi-dot "thousands" price strings) so a plain
# assignment keeps them as text without any extra coercion handling needed.
$ws.Range("D2").Value = '42.056.06'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '2.272.04'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("E12").Value = '  -1.85%  '
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '2.622.86'
$ws.Range("E14").Value = '  +0.79%  '
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '2.276.21'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("E17").Value = '  +4.09%  '
$ws.Range("D18").Value = '41.926.40'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("E19").Value = '  +4.94%  '
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("E24").Value = '  +1.31%  '
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("E29").Value = '  -9.72%  '
$ws.Range("E30").Value = '  +3.64%  '
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("E32").Value = '  +4.23%  '
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("E36").Value = '  +4.20%  '
$ws.Range("E37").Value = '  -1.15%  '
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("E40").Value = '  +0.60%  '
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("E42").Value = '  +2.48%  '
$ws.Range("D43").Value = '2.019.04'
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("E44").Value = '  +8.99%  '
$ws.Range("E45").Value = '  +1.70%  '
$ws.Range("E46").Value = '  +2.71%  '
$ws.Range("E47").Value = '  +2.56%  '
$ws.Range("E48").Value = '  +3.41%  '
$ws.Range("E49").Value = '  +0.71%  '
$ws.Range("E50").Value = '  +3.35%  '
$ws.Range("E51").Value = '  +0.98%  '
